$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily price log gained one new entry. In the sheet's existing (not
# date-sorted) row order, it lands at row 70, pushing the former rows
# 70-184 down to 71-185 (dimension grows from A1:R184 to A1:R185).
$ws.Rows(70).Insert()

# Populate the newly inserted row 70 with the new day's record.
$ws.Cells.Item(70, 1).Value = 6
$ws.Cells.Item(70, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(70, 3).Value = "Metropolitana"
$ws.Cells.Item(70, 4).Value = 44665
$ws.Cells.Item(70, 5).Value = 13
$ws.Cells.Item(70, 6).Value = 100112029
$ws.Cells.Item(70, 7).Value = "Orégano"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 34
$ws.Cells.Item(70, 11).Value = 12000
$ws.Cells.Item(70, 12).Value = 13000
$ws.Cells.Item(70, 13).Value = 12441
$ws.Cells.Item(70, 14).Value = "$/docena de atados"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 4147
$ws.Cells.Item(70, 17).Value = 3
$ws.Cells.Item(70, 18).Value = "Hortaliza"
